$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newly-closed trade as row 9.
$ws.Range("A9").Value = 42654.74627314815
$ws.Range("B9").Value = $false
$ws.Range("C9").Value = 10134.26
$ws.Range("D9").Value = 10141.870000000001
$ws.Range("E9").Value = 308
$ws.Range("F9").Value = 308.45999999999998
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = 0.15
$ws.Range("I9").Value = $false

# Reuse the formatting (date number format etc.) from the previous trade
# row instead of letting Excel mint a brand-new number format for the
# date-typed cells, so row 9 ends up styled exactly like rows 3-8.
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)

# Widen column A so the new (longer) timestamp text still fits, matching
# the repeater's updated "best fit" width.
$ws.Columns.Item(1).ColumnWidth = 14.5
